$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# The old single-scenario "appraisal ratio" test (row 17) is replaced by a
# trio of rows at the bottom of the table covering the three option values
# exercised by the now-completed Appraisal_Ratio macro test
# (option=appraisal / modified / alternative). Removing row 17 shifts every
# following row up by one.
$ws.Rows.Item(17).Delete()

# Append the three new test rows (65 existing rows remain -> new rows 66-68).
# Values are entered column-by-column (all of column A, then B, then C) so
# that the new shared-string entries land in the same order they were
# introduced in the workbook.
$ws.Cells.Item(66, 1).Value2 = "appraisal ratio1"
$ws.Cells.Item(67, 1).Value2 = "appraisal ratio2"
$ws.Cells.Item(68, 1).Value2 = "appraisal ratio3"

$ws.Cells.Item(66, 2).Value2 = "Test appraisal ratio with option=appraisal"
$ws.Cells.Item(67, 2).Value2 = "Test appraisal ratio with option=modified"
$ws.Cells.Item(68, 2).Value2 = "Test appraisal ratio with option=alternative"

$ws.Cells.Item(66, 3).Value2 = "appraisal_ratio_test1"
$ws.Cells.Item(67, 3).Value2 = "appraisal_ratio_test2"
$ws.Cells.Item(68, 3).Value2 = "appraisal_ratio_test3"

# Leave the view scrolled/selected near where the author ended up editing.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 46
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("E55").Select()
